$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N, shifting N:P -> O:Q
$ws.Range("N1:N15").Insert(-4161)

# Make "Repayment Schedule" the active sheet and select full column N
# (matches the recorded selection after the column insert)
$ws.Activate()
$ws.Range("N1:N1048576").Select()
